$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.192.93"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "1.438.88"
$ws.Range("E3").Value = "  +2.59%  "
$ws.Range("E4").Value = "  +1.19%  "
$ws.Range("D5").Value = "'0.9151"
$ws.Range("E5").Value = "  -8.29%  "
$ws.Range("D6").Value = "'275.27"
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("E7").Value = "  -1.41%  "
$ws.Range("D8").Value = "'0.3069"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("D9").Value = "'38.94"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "'1.022"
$ws.Range("E10").Value = "  +2.09%  "
$ws.Range("D11").Value = "'0.06480"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").Value = "'0.9997"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").Value = "'5.334"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "'17.45"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").Value = "'6.043"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").Value = "'0.00001008"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").Value = "1.438.61"
$ws.Range("E17").Value = "  +2.83%  "
$ws.Range("D18").Value = "'0.9339"
$ws.Range("E18").Value = "  -6.28%  "
$ws.Range("D19").Value = "'0.05621"
$ws.Range("E19").Value = "  -1.27%  "
$ws.Range("D20").Value = "'67.57"
$ws.Range("E20").Value = "  -4.61%  "
$ws.Range("D21").Value = "'5.386"
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("D22").Value = "'14.24"
$ws.Range("E22").Value = "  -4.09%  "
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").Value = "'2.238"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").Value = "20.203.95"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("D26").Value = "'138.06"
$ws.Range("E26").Value = "  +1.79%  "
$ws.Range("D27").Value = "'2.123"
$ws.Range("E27").Value = "  -4.85%  "
$ws.Range("D28").Value = "'16.88"
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").Value = "1.588.01"
$ws.Range("E29").Value = "  +2.18%  "
$ws.Range("D30").Value = "'109.82"
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("D31").Value = "'3.826"
$ws.Range("E31").Value = "  -6.54%  "
$ws.Range("D32").Value = "'0.8030"
$ws.Range("E32").Value = "  -2.69%  "
$ws.Range("D33").Value = "'4.821"
$ws.Range("E33").Value = "  -8.69%  "
$ws.Range("D34").Value = "'0.07635"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").Value = "'1.473"
$ws.Range("E35").Value = "  +3.12%  "
$ws.Range("D36").Value = "'0.05838"
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("D37").Value = "'4.652"
$ws.Range("E37").Value = "  -2.81%  "
$ws.Range("D38").Value = "'1.129"
$ws.Range("E38").Value = "  +4.54%  "
$ws.Range("D39").Value = "'0.01981"
$ws.Range("E39").Value = "  -4.49%  "
$ws.Range("D40").Value = "'10.17"
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("D41").Value = "'0.1844"
$ws.Range("E41").Value = "  -2.95%  "
$ws.Range("D42").Value = "'0.9283"
$ws.Range("E42").Value = "  -6.92%  "
$ws.Range("E43").Value = "  -14.23%  "
$ws.Range("D44").Value = "'0.5200"
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("D45").Value = "'3.486"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").Value = "'11.83"
$ws.Range("E46").Value = "  -3.76%  "
$ws.Range("D47").Value = "'116.52"
$ws.Range("E47").Value = "  +4.87%  "
$ws.Range("D48").Value = "'0.5077"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("D50").Value = "'0.06334"
$ws.Range("E50").Value = "  +3.27%  "
$ws.Range("D51").Value = "'0.9880"
$ws.Range("E51").Value = "  -0.68%  "
